# Purchase workbook update — "need to work on fin analysis"
#
# The sheet is a single-row (row 3) monthly rollup table. This edit:
#   - adds a new trailing category "Other Purchases (Trading PL)" (header AU2,
#     data AU3:AW3 = Kgs/Value/Rate, all zero for this month)
#   - fills in the previously-empty CP/Kgs + CP/Value pair (H3:I3) and
#     recalculates the CP/Rate (J3) that depends on it
#   - rolls the whole table's figures forward to a new month (A3: 1-Jan-2025)
#     with refreshed totals across every Kgs/Value/Rate triple
#   - re-homes the window scroll/selection to the newly added columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new category header (shared string) -------------------------------
$ws.Range("AU2").Value = "Other Purchases (Trading PL)"

# ---- row 3 data, column by column (A..AH) -------------------------------
$ws.Range("A3").Value = 45658
$ws.Range("B3").Value = 110000
$ws.Range("C3").Value = 10922850
$ws.Range("D3").Value = 99.298636363636362
$ws.Range("E3").Value = 1574.7
$ws.Range("F3").Value = 369238.8
$ws.Range("G3").Value = 234.48199657077538
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 13140
$ws.Range("J3").Value = 1877.1428571428571
$ws.Range("K3").Value = 111581.7
$ws.Range("L3").Value = 11305228.800000001
$ws.Range("M3").Value = 101.3179472978096
$ws.Range("N3").Value = 740450
$ws.Range("O3").Value = 359214
$ws.Range("P3").Value = 0.48512931325545278
$ws.Range("Q3").Value = 696
$ws.Range("R3").Value = 138061
$ws.Range("S3").Value = 198.36350574712642
$ws.Range("T3").Value = 1629.25
$ws.Range("U3").Value = 335482.5
$ws.Range("V3").Value = 205.91222955347553
$ws.Range("W3").Value = 7455.4
$ws.Range("X3").Value = 1079701
$ws.Range("Y3").Value = 144.82133755398772
$ws.Range("Z3").Value = 9780.65
$ws.Range("AA3").Value = 1553244.5
$ws.Range("AB3").Value = 158.80790131535227
$ws.Range("AC3").Value = 23650.59
$ws.Range("AD3").Value = 5593162.5999999996
$ws.Range("AE3").Value = 236.49146173520404
$ws.Range("AF3").Value = 20175
$ws.Range("AG3").Value = 2432350
$ws.Range("AH3").Value = 120.56257744733581

# AK3 and AN3:AQ3 are unchanged by this edit (left as-is)

# AR3:AT3 mirror the TOTAL Kgs/Value/Rate block (Z3:AB3)
$ws.Range("AR3").Value = 9780.65
$ws.Range("AS3").Value = 1553244.5
$ws.Range("AT3").Value = 158.80790131535227

# ---- new trailing category block (AU3:AW3) -------------------------------
$ws.Range("AU3").Value = 0
$ws.Range("AV3").Value = 0
$ws.Range("AW3").Value = 0

# ---- window state: scroll to the newly added columns & select AJ12 -------
$ws.Range("AJ12").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 27
$win.ScrollRow = 1
